$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: new student record
$ws.Range("A5").Value = 5
$ws.Range("B5").Value = "popescu"
$ws.Range("C5").Value = "ionel"

# D5:F5 hold a numeric-looking grade ("8.0") that must stay stored as TEXT
# (matching the rest of the sheet, e.g. "9.0"/"9.5"/"10.0" in earlier rows).
# Force text entry via NumberFormat "@" so it isn't auto-coerced to a number,
# then clear the formatting again so the cells keep the workbook's default
# (unstyled) look, same as the existing text-number cells in the sheet.
$ws.Range("D5:F5").NumberFormat = "@"
$ws.Range("D5").Value = "8.0"
$ws.Range("E5").Value = "8.0"
$ws.Range("F5").Value = "8.0"
$ws.Range("D5:F5").ClearFormats()

$ws.Range("G5").Value = "budget"
